{"js": "// The document contains several tables with a \"Date\" column whose values\n// are dates in dd.MM.2020 format (23.08.2020 .. 26.08.2020). Every one of\n// those dates needs to be bumped forward by one day (e.g. 23.08.2020 ->\n// 24.08.2020, 26.08.2020 -> 27.08.2020).\n//\n// We do this with a series of body.search()/insertText(\"...\", \"Replace\")\n// passes, one per distinct source date string. To avoid a value that was\n// just written by an earlier pass being matched again by a later pass\n// (e.g. turning 23.08.2020 into 24.08.2020 and then, in a later pass,\n// turning that freshly-written 24.08.2020 into 25.08.2020), the passes run\n// from the latest date to the earliest date.\nconst mapping = [\n  [\"26.08.2020\", \"27.08.2020\"],\n  [\"25.08.2020\", \"26.08.2020\"],\n  [\"24.08.2020\", \"25.08.2020\"],\n  [\"23.08.2020\", \"24.08.2020\"],\n];\n\nfor (const [from, to] of mapping) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains several tables with a \"Date\" column whose values\n# are dates in dd.MM.2020 format (23.08.2020 .. 26.08.2020). Every one of\n# those dates needs to be bumped forward by one day (e.g. 23.08.2020 ->\n# 24.08.2020, 26.08.2020 -> 27.08.2020).\n#\n# We do this with Find/Replace (wdReplaceAll) passes, one per distinct\n# source date string. To avoid a value that was just written by an earlier\n# pass being matched again by a later pass (e.g. turning 23.08.2020 into\n# 24.08.2020 and then, in a later pass, turning that freshly-written\n# 24.08.2020 into 25.08.2020), the passes run from the latest date to the\n# earliest date.\n\n$d = $word.ActiveDocument\n\n$mapping = @(\n  @(\"26.08.2020\", \"27.08.2020\"),\n  @(\"25.08.2020\", \"26.08.2020\"),\n  @(\"24.08.2020\", \"25.08.2020\"),\n  @(\"23.08.2020\", \"24.08.2020\")\n)\n\nforeach ($pair in $mapping) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
